$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price figures formatted as plain text (e.g. "63.821.33",
# "1.00"). Excel auto-detects such strings as numbers on assignment, which would
# mangle formatting (trailing zeros, thousand-dot grouping) and introduce float
# rounding noise. Force the cell to literal text, assign, then clear the
# formatting override so the cell keeps the workbook's original (default) style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.821.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.324.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.74"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.321.64"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.43%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.869.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.323.67"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.854.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.63%  "
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  +5.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.41"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0745"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.88%  "
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "434.46"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.105.07"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.28%  "
$ws.Range("E42").Value = "  +9.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +16.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("E51").Value = "  +3.29%  "
